$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row content: type Name / Level / Description first, then insert an
# "Ancestry" column between Name and Level (matches the shared-string /
# column order captured in the target workbook).
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Level"
$ws.Range("C1").Value = "Description"

$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "Ancestry"

# Column formatting for the four header columns (~15.71 chars wide).
$ws.Columns("A:D").ColumnWidth = 14.8

# Leave the selection on A2, like the saved workbook.
$ws.Range("A2").Select()
